$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first worksheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1395
$wsExpo.Range("F3").Value = 2937

# Sheet "全部类型" (All types) - fourth worksheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1395
$wsAll.Range("F4").Value = 2937
